$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 36
$ws.Range("B6").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B15").Value = 38
$ws.Range("C15").Value = 71
